$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D stores prices as plain text (it keeps thousands-dot grouping and
# trailing zeros, e.g. "64.287.33" / "1.00" / "0.0780"). Whenever the new price
# string looks like a plain Excel number, force the cell to Text format first
# so the COM ".Value =" assignment below does not silently convert it to a
# numeric value and lose formatting / type.

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"

# Update Price (D) and Volume(1h) (E) columns for rows with changed values
$ws.Range("D2").Value = "64.261.25"
$ws.Range("E2").Value = "  +0.13%  "
$ws.Range("D3").Value = "3.491.44"
$ws.Range("E3").Value = "  -1.00%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "587.42"
$ws.Range("E5").Value = "  +0.26%  "
$ws.Range("D6").Value = "134.14"
$ws.Range("E6").Value = "  +0.68%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("E9").Value = "  -0.34%  "
$ws.Range("E10").Value = "  +2.04%  "
$ws.Range("D11").Value = "0.385"
$ws.Range("E11").Value = "  +1.74%  "
$ws.Range("D12").Value = "4.087.05"
$ws.Range("E12").Value = "  -1.14%  "
$ws.Range("E13").Value = "  +1.02%  "
$ws.Range("E14").Value = "  +1.29%  "
$ws.Range("D15").Value = "3.495.11"
$ws.Range("E15").Value = "  -0.32%  "
$ws.Range("D16").Value = "25.79"
$ws.Range("E16").Value = "  -6.60%  "
$ws.Range("D17").Value = "64.348.92"
$ws.Range("E17").Value = "  +0.14%  "
$ws.Range("D18").Value = "9.86"
$ws.Range("E18").Value = "  +0.95%  "
$ws.Range("E19").Value = "  +2.27%  "
$ws.Range("D20").Value = "13.60"
$ws.Range("E20").Value = "  -3.20%  "
$ws.Range("D21").Value = "393.88"
$ws.Range("E21").Value = "  +2.20%  "
$ws.Range("E22").Value = "  -0.80%  "
$ws.Range("D23").Value = "3.631.74"
$ws.Range("E23").Value = "  -1.09%  "
$ws.Range("E24").Value = "  +0.92%  "
$ws.Range("E25").Value = "  -0.01%  "
$ws.Range("D26").Value = "5.74"
$ws.Range("E26").Value = "  +0.14%  "
$ws.Range("E27").Value = "  -0.38%  "
$ws.Range("E30").Value = "  -0.01%  "
$ws.Range("E31").Value = "  -4.82%  "
$ws.Range("E32").Value = "  -2.11%  "
$ws.Range("D33").Value = "3.514.09"
$ws.Range("E33").Value = "  -0.77%  "
$ws.Range("E34").Value = "  +3.42%  "
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("D36").Value = "23.39"
$ws.Range("E36").Value = "  -0.91%  "
$ws.Range("E37").Value = "  -4.18%  "
$ws.Range("D38").Value = "6.88"
$ws.Range("E38").Value = "  -0.46%  "
$ws.Range("E39").Value = "  -0.63%  "
$ws.Range("D40").Value = "166.63"
$ws.Range("E40").Value = "  +3.52%  "
$ws.Range("D41").Value = "0.0778"
$ws.Range("E41").Value = "  -1.26%  "
$ws.Range("E42").Value = "  -1.21%  "
$ws.Range("E43").Value = "  -0.09%  "
$ws.Range("D44").Value = "25.27"
$ws.Range("E44").Value = "  -4.78%  "
$ws.Range("D45").Value = "4.38"
$ws.Range("E45").Value = "  -0.70%  "
$ws.Range("E46").Value = "  +2.04%  "
$ws.Range("E47").Value = "  -3.70%  "
$ws.Range("D48").Value = "2.459.38"
$ws.Range("D49").Value = "6.75"
$ws.Range("E49").Value = "  -0.90%  "
$ws.Range("E50").Value = "  -1.44%  "
$ws.Range("E51").Value = "  -1.20%  "

# Rows 28 and 29 had their content swapped (RenderToken <-> Binance-PegBSC-USD)
$ws.Range("B28").Value = "Binance-PegBSC-USD"
$ws.Range("C28").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D28").Value = "0.995"
$ws.Range("E28").Value = "  -0.45%  "

$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").Value = "7.36"
$ws.Range("E29").Value = "  -1.64%  "
